$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update totals
$ws.Range("E11").Value = 9280
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Remove the second worker row (row 17) entirely, shifting rows 18+ up by one
$ws.Rows("17").Delete()
